$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author deleted the 10 data rows covering 26/03/2025 - 29/03/2025
# (rows 1398-1407), which shifts everything below them up by 10 rows and
# drops the now-unused "26/03/2025".."29/03/2025" shared-string entries.
$ws.Range("A1398:G1407").EntireRow.Delete()

# Leave the view parked on the row that was being edited, matching where
# the author's selection ended up after the delete.
$ws.Range("A1397").Select()
